# Replace the 100 arithmetic expressions in the table, cell by cell,
# in document (row-major) order, since several expressions repeat
# (e.g. "52-17=" occurs twice with different replacements).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$oldVals = @('98-49=', '17+47=', '27+29=', '88+7=', '45+29=', '63-45=', '19+12=', '92-33=', '96-9=', '47+28=', '26+48=', '92-56=', '79+13=', '92-58=', '70-67=', '71-12=', '50-18=', '48+45=', '77+8=', '62-14=', '39+42=', '90-65=', '29+62=', '95-19=', '44+47=', '21-9=', '85-37=', '83-8=', '31-9=', '83+8=', '58+35=', '92-6=', '9+63=', '31-27=', '22-5=', '8+58=', '24+17=', '32-17=', '97-18=', '94-38=', '82+9=', '51-16=', '18+59=', '91-54=', '24+28=', '19+28=', '90-16=', '80-44=', '96-29=', '36-7=', '61-37=', '54-9=', '34-17=', '52-17=', '25-7=', '40-6=', '8+84=', '93-48=', '18+66=', '92-3=', '81-62=', '98-69=', '27+18=', '80-46=', '81-44=', '26+45=', '84-48=', '14+67=', '8+24=', '24+29=', '75+19=', '95-86=', '41-17=', '66-57=', '91-74=', '95-77=', '45-19=', '82-35=', '77+5=', '57+9=', '5+77=', '52-17=', '35+46=', '17+15=', '45+9=', '18+43=', '42-14=', '88+9=', '80-9=', '91-26=', '9+14=', '90-35=', '59+29=', '58+34=', '50-35=', '94-9=', '63-37=', '18+6=', '93-86=', '86+9=')
$newVals = @('34+18=', '77-39=', '62-49=', '18+3=', '83-55=', '62-46=', '39+32=', '38+57=', '39+54=', '40-4=', '4+69=', '40-39=', '57-48=', '92-66=', '45+16=', '54+19=', '16+38=', '15+38=', '46+5=', '17+25=', '52-3=', '54+18=', '92-27=', '7+17=', '55-17=', '9+24=', '80-21=', '22+9=', '60-15=', '13+78=', '38+59=', '42-8=', '39+36=', '69+3=', '19+22=', '56-38=', '69+27=', '8+78=', '46+29=', '48+44=', '39+49=', '28+68=', '2+79=', '28+24=', '93-39=', '38+24=', '56-27=', '94-55=', '72+19=', '32+29=', '38+15=', '94-89=', '84-8=', '61-5=', '50-9=', '55-9=', '16+69=', '59+37=', '37+58=', '24+67=', '32+29=', '16+25=', '87-9=', '95-87=', '53+29=', '90-75=', '45-29=', '7+39=', '23+38=', '37+37=', '25+48=', '39+6=', '79+15=', '30-22=', '25+46=', '82-3=', '93-39=', '55-29=', '65-49=', '29+15=', '39+33=', '53-15=', '89+7=', '52-33=', '29+67=', '9+64=', '74-59=', '22-19=', '4+19=', '61-16=', '75-69=', '48-9=', '39+37=', '77+15=', '23+29=', '16+6=', '62+19=', '18+3=', '12-9=', '7+18=')

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $expected = $oldVals[$idx]
        $replacement = $newVals[$idx]
        $current = $cell.Range.Text.Substring(0, $expected.Length)
        if ($current -ne $expected) {
            Write-Host "MISMATCH at idx" $idx "row" $r "col" $c "expected" $expected "got" $current
        }
        $cell.Range.Text = $replacement
        $idx = $idx + 1
    }
}
Write-Host "Done. Processed" $idx "cells."
